# Update loading_percent results for the 380 kV case (row-wise per-column values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.14795727283463
$ws.Range("C2").Value = 10.55635688124185
$ws.Range("E2").Value = 15.63813085493678
$ws.Range("F2").Value = 41.39517530451591
$ws.Range("G2").Value = 3.684834061996117
$ws.Range("J2").Value = 8.76890932344336
$ws.Range("K2").Value = 9.782609893465249
$ws.Range("L2").Value = 11.99414111826208
$ws.Range("O2").Value = 28.05351144098294
$ws.Range("B3").Value = 13.92443809128069
$ws.Range("C3").Value = 10.57707018133788
$ws.Range("E3").Value = 15.63398402917329
$ws.Range("F3").Value = 41.43800736263194
$ws.Range("G3").Value = 3.686698010628028
$ws.Range("J3").Value = 8.771442268996324
$ws.Range("K3").Value = 9.624193686114985
$ws.Range("L3").Value = 11.97088608352323
$ws.Range("O3").Value = 28.15510614397759
$ws.Range("B4").Value = 13.7877875818843
$ws.Range("C4").Value = 10.59049630267863
$ws.Range("E4").Value = 15.63380293202275
$ws.Range("F4").Value = 41.47340978334444
$ws.Range("G4").Value = 3.687902996611994
$ws.Range("J4").Value = 8.773191473528511
$ws.Range("K4").Value = 9.526977778973251
$ws.Range("L4").Value = 11.95824025331951
$ws.Range("O4").Value = 28.22273393862956
$ws.Range("B5").Value = 13.73231912629837
$ws.Range("C5").Value = 10.59614611337222
$ws.Range("E5").Value = 15.63432569612871
$ws.Range("F5").Value = 41.49012386683177
$ws.Range("G5").Value = 3.688409302512838
$ws.Range("J5").Value = 8.773953243877965
$ws.Range("K5").Value = 9.487422472546225
$ws.Range("L5").Value = 11.95350119902789
$ws.Range("O5").Value = 28.25161106061371
$ws.Range("B6").Value = 13.72312380698254
$ws.Range("C6").Value = 10.59709505950918
$ws.Range("E6").Value = 15.63444858095169
$ws.Range("F6").Value = 41.49303730604679
$ws.Range("G6").Value = 3.68849429755015
$ws.Range("J6").Value = 8.774082697353744
$ws.Range("K6").Value = 9.480859432686268
$ws.Range("L6").Value = 11.95273940152951
$ws.Range("O6").Value = 28.25648566060707
$ws.Range("B7").Value = 13.78703854178872
$ws.Range("C7").Value = 10.59057177435242
$ws.Range("E7").Value = 15.63380756447358
$ws.Range("F7").Value = 41.47362593712874
$ws.Range("G7").Value = 3.687909762962282
$ws.Range("J7").Value = 8.773201548574541
$ws.Range("K7").Value = 9.526444011911877
$ws.Range("L7").Value = 11.95817465900166
$ws.Range("O7").Value = 28.22311805002837
$ws.Range("B8").Value = 14.07080647129767
$ws.Range("C8").Value = 10.56335217192042
$ws.Range("E8").Value = 15.63621133298689
$ws.Range("F8").Value = 41.40805355452097
$ws.Range("G8").Value = 3.685464220507692
$ws.Range("J8").Value = 8.769742533823802
$ws.Range("K8").Value = 9.72800565022791
$ws.Range("L8").Value = 11.98578613787968
$ws.Range("O8").Value = 28.08745114840694
$ws.Range("B9").Value = 14.62896171652516
$ws.Range("C9").Value = 10.51556984027779
$ws.Range("E9").Value = 15.65960162124306
$ws.Range("F9").Value = 41.35174580608121
$ws.Range("G9").Value = 3.68114649602339
$ws.Range("J9").Value = 8.764490595428493
$ws.Range("K9").Value = 10.12161570077297
$ws.Range("L9").Value = 12.05271846411261
$ws.Range("O9").Value = 27.86311611421303
$ws.Range("B10").Value = 15.03612494223855
$ws.Range("C10").Value = 10.48384314828742
$ws.Range("E10").Value = 15.68804038740247
$ws.Range("F10").Value = 41.35445818437199
$ws.Range("G10").Value = 3.678262603623592
$ws.Range("J10").Value = 8.761555035733771
$ws.Range("K10").Value = 10.40709759809119
$ws.Range("L10").Value = 12.10944644472049
$ws.Range("O10").Value = 27.72381110266386
$ws.Range("B11").Value = 15.21991995352962
$ws.Range("C11").Value = 10.47013681163775
$ws.Range("E11").Value = 15.70338855199025
$ws.Range("F11").Value = 41.3652488361864
$ws.Range("G11").Value = 3.677012608184878
$ws.Range("J11").Value = 8.760417804766165
$ws.Range("K11").Value = 10.53562232591092
$ws.Range("L11").Value = 12.13683660476094
$ws.Range("O11").Value = 27.66599652229385
$ws.Range("B12").Value = 15.28924950428938
$ws.Range("C12").Value = 10.4650505111757
$ws.Range("E12").Value = 15.70954398029002
$ws.Range("F12").Value = 41.37070616780417
$ws.Range("G12").Value = 3.676548119792136
$ws.Range("J12").Value = 8.76001547834853
$ws.Range("K12").Value = 10.58405512770397
$ws.Range("L12").Value = 12.14743111715379
$ws.Range("O12").Value = 27.64490419410849
$ws.Range("B13").Value = 15.27433113729105
$ws.Range("C13").Value = 10.46614131900433
$ws.Range("E13").Value = 15.70820308089404
$ws.Range("F13").Value = 41.36946991196854
$ws.Range("G13").Value = 3.676647762417951
$ws.Range("J13").Value = 8.760100869916487
$ws.Range("K13").Value = 10.57363546013832
$ws.Range("L13").Value = 12.14513959372767
$ws.Range("O13").Value = 27.64941115928832
$ws.Range("B14").Value = 15.22562950181347
$ws.Range("C14").Value = 10.46971627703641
$ws.Range("E14").Value = 15.70388810021375
$ws.Range("F14").Value = 41.36567035573568
$ws.Range("G14").Value = 3.676974217171219
$ws.Range("J14").Value = 8.760384138610384
$ws.Range("K14").Value = 10.53961190915945
$ws.Range("L14").Value = 12.13770379498553
$ws.Range("O14").Value = 27.6642451840002
$ws.Range("B15").Value = 15.19576135763708
$ws.Range("C15").Value = 10.47191957164144
$ws.Range("E15").Value = 15.70128966943733
$ws.Range("F15").Value = 41.36352147200401
$ws.Range("G15").Value = 3.677175332385971
$ws.Range("J15").Value = 8.760561331448146
$ws.Range("K15").Value = 10.51873941048346
$ws.Range("L15").Value = 12.13317796006744
$ws.Range("O15").Value = 27.6734357948418
$ws.Range("B16").Value = 15.02407923381707
$ws.Range("C16").Value = 10.48475346605182
$ws.Range("E16").Value = 15.68708559075819
$ws.Range("F16").Value = 41.35394505421902
$ws.Range("G16").Value = 3.678345535646931
$ws.Range("J16").Value = 8.761633328891282
$ws.Range("K16").Value = 10.39866741317885
$ws.Range("L16").Value = 12.10768786037348
$ws.Range("O16").Value = 27.72770138096946
$ws.Range("B17").Value = 14.91834586353019
$ws.Range("C17").Value = 10.49281234506189
$ws.Range("E17").Value = 15.67898713965982
$ws.Range("F17").Value = 41.35051606263552
$ws.Range("G17").Value = 3.679079240790593
$ws.Range("J17").Value = 8.76234159375826
$ws.Range("K17").Value = 10.32463161150768
$ws.Range("L17").Value = 12.09245254828208
$ws.Range("O17").Value = 27.76241590499921
$ws.Range("B18").Value = 14.85740034907772
$ws.Range("C18").Value = 10.49751598785455
$ws.Range("E18").Value = 15.67455629565666
$ws.Range("F18").Value = 41.3494435763442
$ws.Range("G18").Value = 3.679507077648061
$ws.Range("J18").Value = 8.762767632307574
$ws.Range("K18").Value = 10.28192431621746
$ws.Range("L18").Value = 12.08383906019041
$ws.Range("O18").Value = 27.78290562513923
$ws.Range("B19").Value = 14.83674484835437
$ws.Range("C19").Value = 10.49912032007887
$ws.Range("E19").Value = 15.6730952018494
$ws.Range("F19").Value = 41.3492350669416
$ws.Range("G19").Value = 3.679652938340693
$ws.Range("J19").Value = 8.762915092522858
$ws.Range("K19").Value = 10.26744444955728
$ws.Range("L19").Value = 12.08094851073103
$ws.Range("O19").Value = 27.78993283440256
$ws.Range("B20").Value = 14.92961534178062
$ws.Range("C20").Value = 10.49194738892463
$ws.Range("E20").Value = 15.67982574552652
$ws.Range("F20").Value = 41.35078798016928
$ws.Range("G20").Value = 3.679000533656385
$ws.Range("J20").Value = 8.762264267458796
$ws.Range("K20").Value = 10.33252599070481
$ws.Range("L20").Value = 12.09405894279891
$ws.Range("O20").Value = 27.75866635871119
$ws.Range("B21").Value = 15.23994217210802
$ws.Range("C21").Value = 10.46866340644314
$ws.Range("E21").Value = 15.70514622225374
$ws.Range("F21").Value = 41.36674919478207
$ws.Range("G21").Value = 3.676878089446418
$ws.Range("J21").Value = 8.760300168599326
$ws.Range("K21").Value = 10.54961221379119
$ws.Range("L21").Value = 12.13988187598698
$ws.Range("O21").Value = 27.65986632321766
$ws.Range("B22").Value = 15.4411590505419
$ws.Range("C22").Value = 10.45405191748496
$ws.Range("E22").Value = 15.72369456558406
$ws.Range("F22").Value = 41.38517070115892
$ws.Range("G22").Value = 3.675542559214073
$ws.Range("J22").Value = 8.759181494309932
$ws.Range("K22").Value = 10.69009196030525
$ws.Range("L22").Value = 12.1711237821503
$ws.Range("O22").Value = 27.59996330219246
$ws.Range("B23").Value = 15.33393282769697
$ws.Range("C23").Value = 10.46179504728504
$ws.Range("E23").Value = 15.71361313189267
$ws.Range("F23").Value = 41.37460900603189
$ws.Range("G23").Value = 3.676250648864136
$ws.Range("J23").Value = 8.759763515246075
$ws.Range("K23").Value = 10.61525716881467
$ws.Range("L23").Value = 12.15433282571813
$ws.Range("O23").Value = 27.63150689461623
$ws.Range("B24").Value = 14.92452089781873
$ws.Range("C24").Value = 10.49233821600422
$ws.Range("E24").Value = 15.67944591052535
$ws.Range("F24").Value = 41.35066224604888
$ws.Range("G24").Value = 3.679036098396571
$ws.Range("J24").Value = 8.762299167951314
$ws.Range("K24").Value = 10.32895738432959
$ws.Range("L24").Value = 12.09333223823598
$ws.Range("O24").Value = 27.76035987193002
$ws.Range("B25").Value = 14.47819839485066
$ws.Range("C25").Value = 10.5279005507609
$ws.Range("E25").Value = 15.65128708814093
$ws.Range("F25").Value = 41.3592349115696
$ws.Range("G25").Value = 3.682263700534505
$ws.Range("J25").Value = 8.765748536840237
$ws.Range("K25").Value = 10.01559394146862
$ws.Range("L25").Value = 12.03326656059662
$ws.Range("O25").Value = 27.91933049282319
